$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Change 1: "[生产许可证号] " -> split into three runs:
#   "[生产许可证"  +  "编"  +  "号] "
# (same bold "宋体" run formatting throughout; a new character "编" is
#  inserted between "证" and "号").
# -----------------------------------------------------------------------
$labelScope = $d.Content
$labelScope.Find.Execute("[生产许可证号] ")
if ($labelScope.Find.Found) {
    $labelStart = $labelScope.Start
    $labelEnd = $labelScope.End

    # Locate "号] " inside that run so we know exactly where to split /
    # insert the new character, regardless of absolute offsets.
    $tailScope = $d.Range($labelStart, $labelEnd)
    $tailScope.Find.Execute("号] ")
    $splitPos = $tailScope.Start

    # Insert the new character "编" right before "号] ".
    $insertPoint = $d.Range($splitPos, $splitPos)
    $insertPoint.InsertBefore("编")

    # Isolate the freshly inserted "编" into its own run by toggling a
    # character property off and back on - this forces the engine to
    # break it out of the surrounding run instead of re-absorbing it.
    $newCharRange = $d.Range($splitPos, $splitPos + 1)
    $newCharRange.Bold = 9999999
    $newCharRange.Bold = 1
}

# -----------------------------------------------------------------------
# Change 2: "1721 " -> "1722 " (energy value), runs must stay split from
# the following "千焦" run.
# -----------------------------------------------------------------------
$energyScope = $d.Content
$energyScope.Find.Execute("1721")
if ($energyScope.Find.Found) {
    $energyRange = $d.Range($energyScope.Start, $energyScope.End)
    $energyRange.Text = "1722"

    # The text replace above can fuse the "1722 " run back together with
    # the adjacent "千焦" run (identical formatting). Re-split them so the
    # two runs remain distinct, matching the original structure.
    $afterScope = $d.Range($energyScope.Start, $energyScope.Start + 40)
    $afterScope.Find.Execute("千焦")
    if ($afterScope.Find.Found) {
        $qiRange = $d.Range($afterScope.Start, $afterScope.End)
        $qiRange.Bold = 9999999
        $qiRange.Bold = 1
    }
}

# -----------------------------------------------------------------------
# Change 3: "20%" -> "21%" (NRV% for energy). This run is alone in its
# cell/paragraph so a plain replace is safe.
# -----------------------------------------------------------------------
$pctScope = $d.Content
$pctScope.Find.Execute("20%", $false, $false, $false, $false, $false, $true, 1, $false, "21%", 2)
